# Apply "repull data, push all data, mean calculation" update:
# Updates column F (dSF) values for rows 2-24 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = -5
    3  = 1
    4  = -1
    5  = -8
    6  = 4
    7  = -2
    8  = -1
    9  = 1
    10 = -3
    11 = -1
    12 = -3
    13 = 2
    14 = -2
    15 = 2
    16 = 5
    17 = -5
    18 = -1
    19 = 4
    20 = 3
    21 = 4
    22 = -2
    23 = -3
    24 = 0
}

foreach ($row in $newValues.Keys) {
    $ws.Range("F$row").Value = $newValues[$row]
}
